# partbform.xlsx - "Add files via upload"
#
# The header row is simplified: the long descriptive headers (one of them
# rich-text with a superscript "-1" exponent) are replaced with short plain
# labels, the now-unneeded extra row height collapses, and the active
# selection moves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: "Glucose concentration (mM)" -> "glucose"
$ws.Range("A1").Value = "glucose"

# B1: rich-text "Growth rate (h-1)" (with superscript "-1") -> plain "growthrate"
# Assigning a plain string clears the old per-run (bold/superscript) formatting,
# leaving the cell's own font (bold Arial) as the only formatting.
$ws.Range("B1").Value = "growthrate"

# Row 1 no longer needs to wrap the long two-part header text, so it shrinks.
$ws.Rows.Item(1).RowHeight = 35

# Move the active selection to E8.
$ws.Range("E8").Select() | Out-Null
